$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3046
$ws.Range("B2").Value = "Sra. Giovanna Castro"
$ws.Range("C2").Value = "TI"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45095
$ws.Range("G2").Value = 11978.02

$ws.Range("A3").Value = 30504
$ws.Range("B3").Value = "Dra. Eloah Viana"
$ws.Range("C3").Value = "Financeiro"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5150.23

$ws.Range("A4").Value = 99205
$ws.Range("B4").Value = "Srta. Lavínia Ramos"
$ws.Range("D4").Value = "Viagem de negócios"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45093
$ws.Range("G4").Value = 8892.17

$ws.Range("A5").Value = 50990
$ws.Range("B5").Value = "Diogo Rezende"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("F5").Value = 45090
$ws.Range("G5").Value = 6406.24

$ws.Range("A6").Value = 24283
$ws.Range("B6").Value = "Ana Carolina Pinto"
$ws.Range("C6").Value = "TI"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("F6").Value = 45084
$ws.Range("G6").Value = 5209.7

$ws.Range("A7").Value = 71456
$ws.Range("B7").Value = "Alice Fogaça"
$ws.Range("C7").Value = "Vendas"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45083
$ws.Range("G7").Value = 9332.02

$ws.Range("A8").Value = 80752
$ws.Range("B8").Value = "Dra. Isadora Ferreira"
$ws.Range("C8").Value = "P&D"
$ws.Range("D8").Value = "Doença"
$ws.Range("E8").Value = 6
$ws.Range("G8").Value = 4179.26

$ws.Range("A9").Value = 93468
$ws.Range("B9").Value = "João Felipe Araújo"
$ws.Range("C9").Value = "Jurídico"
$ws.Range("D9").Value = "Viagem de negócios"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45081
$ws.Range("G9").Value = 3177.99

$ws.Range("A10").Value = 2386
$ws.Range("B10").Value = "Nicolas Cavalcanti"
$ws.Range("C10").Value = "Vendas"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45082
$ws.Range("G10").Value = 7282.65

$ws.Range("A11").Value = 59238
$ws.Range("B11").Value = "Laura Castro"
$ws.Range("F11").Value = 45080
$ws.Range("G11").Value = 4364.22
